# Apply BoM reference-designator corrections to the "BoM" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

$ws.Range("D9").Value  = "C6"
$ws.Range("D11").Value = "C5 C8 C11 C12 C13"
$ws.Range("D12").Value = "C1 C2 C7 C10"
$ws.Range("D18").Value = "R3 R4"
$ws.Range("D19").Value = "R2 R5 R6 R7 R8 R9 R10 R11"
$ws.Range("D21").Value = "U1"
$ws.Range("D22").Value = "U2"
